$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4.55194436134888
$ws.Range("E2").Value = 0.0563609197510995
$ws.Range("F2").Value = 9.45954957352126
$ws.Range("G2").Value = 1
$ws.Range("D3").Value = -1.648950251367
$ws.Range("E3").Value = -13.5052864425289
$ws.Range("F3").Value = 11.2187000547253
$ws.Range("G3").Value = 0
$ws.Range("D4").Value = -6.97982475932386
$ws.Range("E4").Value = -16.6459358443617
$ws.Range("F4").Value = 4.74889646073614
$ws.Range("G4").Value = 0
$ws.Range("D5").Value = -1.44699984262765
$ws.Range("E5").Value = -19.6318810745464
$ws.Range("F5").Value = 21.0400122736455
$ws.Range("G5").Value = 0
$ws.Range("D6").Value = -45.9166131821395
$ws.Range("E6").Value = -52.1196125977235
$ws.Range("F6").Value = -38.2501219701516
$ws.Range("G6").Value = 1
$ws.Range("D7").Value = 20.3167387785498
$ws.Range("E7").Value = 13.2842041030328
$ws.Range("F7").Value = 28.686158084496
$ws.Range("G7").Value = 1
$ws.Range("D8").Value = 21.3210454396107
$ws.Range("E8").Value = 12.0020624217291
$ws.Range("F8").Value = 31.4976304702006
$ws.Range("G8").Value = 1
$ws.Range("D9").Value = 57.6830514671484
$ws.Range("E9").Value = 31.8625096642734
$ws.Range("F9").Value = 91.9316233452192
$ws.Range("G9").Value = 1
$ws.Range("D10").Value = 56.8380396682825
$ws.Range("E10").Value = 29.6456058789541
$ws.Range("F10").Value = 90.2060849455851
$ws.Range("G10").Value = 1
$ws.Range("D11").Value = -16.5365931509972
$ws.Range("E11").Value = -25.1675826809712
$ws.Range("F11").Value = -5.90321396800044
$ws.Range("G11").Value = 1
$ws.Range("D12").Value = 3.39013510432004
$ws.Range("E12").Value = -0.692014453713557
$ws.Range("F12").Value = 7.89734044144451
$ws.Range("G12").Value = 0
$ws.Range("D13").Value = 19.8336351311463
$ws.Range("E13").Value = 7.29071820799665
$ws.Range("F13").Value = 33.9609753181257
$ws.Range("G13").Value = 1
$ws.Range("D14").Value = 122.247679917388
$ws.Range("E14").Value = 65.5287541530487
$ws.Range("F14").Value = 223.377357672785
$ws.Range("G14").Value = 1
$ws.Range("D15").Value = 204.413127015303
$ws.Range("E15").Value = 106.907756757276
$ws.Range("F15").Value = 434.060730444603
$ws.Range("G15").Value = 1
$ws.Range("D16").Value = 147.522904825849
$ws.Range("E16").Value = 85.0514608735944
$ws.Range("F16").Value = 233.826154721825
$ws.Range("G16").Value = 1
$ws.Range("D17").Value = 3.98639749452161
$ws.Range("E17").Value = 0.69864422689744
$ws.Range("F17").Value = 7.54413189225096
$ws.Range("G17").Value = 1
$ws.Range("D18").Value = 21.785328233901
$ws.Range("E18").Value = 13.516708171832
$ws.Range("F18").Value = 31.8467795250282
$ws.Range("G18").Value = 1
$ws.Range("D19").Value = 76.0688266715239
$ws.Range("E19").Value = 52.254380206788
$ws.Range("F19").Value = 104.80959814208
$ws.Range("G19").Value = 1
$ws.Range("D20").Value = 67.285325009361
$ws.Range("E20").Value = 30.5152582497989
$ws.Range("F20").Value = 116.452643513778
$ws.Range("G20").Value = 1
$ws.Range("D21").Value = -43.3778103466629
$ws.Range("E21").Value = -48.7866056217417
$ws.Range("F21").Value = -37.3549402776375
$ws.Range("G21").Value = 1
